# The "EP1 Neat Resin Neat Resin" entry (A31) was removed from the
# "train material" list in column A. Every subsequent value in column A
# (A32:A133) shifts up by one row; columns B and C are untouched. The
# now-unused last row (133) is cleared so the sheet's used range shrinks
# from A1:C133 to A1:C132.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 133

for ($r = 31; $r -le ($lastRow - 1); $r++) {
    $nextValue = $ws.Cells.Item($r + 1, 1).Value()
    $ws.Cells.Item($r, 1).Value = $nextValue
}

# The old last row's data has all been pulled up one slot; clear what is
# now a trailing empty row so it drops out of the sheet's used range.
$ws.Range("A133:C133").ClearContents()
